# is takip guncellemesi - 12.11.2025 19:04:14
# Tum tarih hucrelerini 1 gun geri kaydirir (J/K sutunlari - Is Takip Listesi,
# cesitli tarih sutunlari - Guncelleme) ve birkac NOTLAR hucresine yeni satir ekler.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    # Hucreye duz metin yazar; tarih/sayi gibi gorunen degerlerin Excel
    # tarafindan otomatik olarak sayiya/taruhe donusturulmesini engellemek
    # icin basina kesme isareti ekler, sonra da bu yuzden olusan
    # "metin bicimi" hucre stilini temizler (stil degismeden kalsin diye).
    param($range, [string]$value)
    $range.Value = "'" + $value
    $range.ClearFormats()
}

function Shift-DateString {
    # "yyyy-MM-dd" bicimli metni 1 gun geriye alir.
    param([string]$s)
    $d = [DateTime]::ParseExact($s, "yyyy-MM-dd", $null)
    $d2 = $d.AddDays(-1)
    $result = $d2.ToString("yyyy-MM-dd")
    return $result
}

function Shift-CellIfDate {
    # Hucre doluysa, icindeki tarihi 1 gun geri al ve metin olarak geri yaz.
    param($cell)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $newval = Shift-DateString $val
        Set-TextValue $cell $newval
    }
}

# ---------------------------------------------------------------------
# Sayfa 1: "İş Takip Listesi" - J (İŞE BAŞLAMA/YER TESLİMİ) ve
# K (İHALE BİTİŞ TARİHİ) sutunlari, tarih iceren her satir icin.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("İş Takip Listesi")

for ($r = 2; $r -le 122; $r++) {
    $jcell = $ws1.Cells.Item($r, 10)
    $kcell = $ws1.Cells.Item($r, 11)
    Shift-CellIfDate $jcell
    Shift-CellIfDate $kcell
}

# NOTLAR (M) sutununda birkac hucreye yeni guncelleme notu ekleniyor.
function Append-Note {
    param($cell, [string]$extra)
    $old = $cell.Value2
    $new = $old + "`n" + $extra
    $cell.Value = $new
}

Append-Note $ws1.Cells.Item(98, 13) "GM OLUR'u beklenecek"
Append-Note $ws1.Cells.Item(100, 13) "10.11.2025 Kroki kontrol için 11.11.2025 de verilecek. 14.11.2025 de firmaya teslim edilecek"
Append-Note $ws1.Cells.Item(102, 13) "10.11.2025 Ormancı 11.11.2025 teslim edecek"
Append-Note $ws1.Cells.Item(107, 13) "10.11.2025 Tutanak 12.11.2025 de basılacak"
Append-Note $ws1.Cells.Item(118, 13) "10.11.2025 Değerlendrime 13.11.2025 de bitecek"
Append-Note $ws1.Cells.Item(121, 13) "10.11.2025 Tutanaklar 13.11.2025 imzaya götürülecek"
Append-Note $ws1.Cells.Item(122, 13) "10.11.2025 12.11.2025 askıya çıkacak(Tutanak kontrolü 15.11.2025 de bitip 17.11.2025 askıya çıkacak)"

# ---------------------------------------------------------------------
# Sayfa 2: "Güncelleme" - F,G,H,I,J,N,P gibi tarih sutunlarinda, her satirda
# degisen hucreler farkli oldugu icin hangi satirda hangi sutunlarin
# degisecegi bir tablo (hashtable) ile belirtiliyor.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Güncelleme")

$colIndex = @{ "F" = 6; "G" = 7; "H" = 8; "I" = 9; "J" = 10; "K" = 11; "N" = 14; "P" = 16 }

$rowChanges = @{
    2  = @("J", "N", "P")
    3  = @("J", "N")
    4  = @("J", "N", "P")
    5  = @("I")
    6  = @("J", "N")
    7  = @("I", "J")
    8  = @("J", "N", "P")
    9  = @("I", "J")
    10 = @("J", "N")
    11 = @("I", "J", "N")
    12 = @("J", "N")
    13 = @("J")
    14 = @("J")
    15 = @("J", "N")
    16 = @("J", "N", "P")
    17 = @("J")
    18 = @("J")
    19 = @("I", "J", "N")
    20 = @("J")
    21 = @("J")
    22 = @("J")
    23 = @("J")
    24 = @("I")
    25 = @("J")
    27 = @("J")
    28 = @("J")
    29 = @("I", "J")
}

foreach ($r in $rowChanges.Keys) {
    $cols = $rowChanges[$r]
    foreach ($colLetter in $cols) {
        $c = $colIndex[$colLetter]
        $cell = $ws2.Cells.Item($r, $c)
        Shift-CellIfDate $cell
    }
}
